# Apply coin price/volume/hour updates per commit "Updated symbol list on Fri Feb 10 11:11:02 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text is a plain label/URL - no special number formatting required.
$textUpdates = @{
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B16" = "UpBots"
    "C16" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B18" = "BTSEToken"
    "C18" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B19" = "BitpandaEcosystemToken"
    "C19" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "B20" = "MCDex"
    "C20" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "B21" = "ProBitToken"
    "C21" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "B22" = "ZBToken"
    "C22" = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
}

# Cells whose text looks numeric (prices, percentages, hour codes) - force Text format
# so Excel keeps the literal string (with trailing zeros, leading minus, % sign, etc.)
# instead of silently coercing it to a Number/Percentage value.
$numericLookingUpdates = @{
    "D2" = "306.41"
    "E2" = "-4.50%"
    "G2" = "11"
    "D3" = "40.02"
    "E3" = "-5.96%"
    "G3" = "11"
    "D4" = "5.124"
    "E4" = "-1.22%"
    "G4" = "11"
    "D5" = "0.07735"
    "E5" = "-5.39%"
    "G5" = "11"
    "D6" = "4.236"
    "E6" = "-1.72%"
    "G6" = "11"
    "D7" = "1.604"
    "E7" = "-11.51%"
    "G7" = "11"
    "D8" = "0.8890"
    "E8" = "-4.83%"
    "G8" = "11"
    "D9" = "0.1000"
    "E9" = "-10.02%"
    "G9" = "11"
    "D10" = "0.1737"
    "E10" = "-7.13%"
    "G10" = "11"
    "D11" = "0.08979"
    "E11" = "-5.04%"
    "G11" = "11"
    "D12" = "0.04457"
    "E12" = "-4.84%"
    "G12" = "11"
    "E13" = "-0.21%"
    "G13" = "11"
    "D14" = "0.001262"
    "E14" = "-2.94%"
    "G14" = "11"
    "D15" = "0.005829"
    "E15" = "2.07%"
    "G15" = "11"
    "D16" = "0.007491"
    "E16" = "2,412.92%"
    "G16" = "11"
    "D17" = "3.353"
    "E17" = "-0.07%"
    "G17" = "11"
    "D18" = "2.410"
    "E18" = "-4.83%"
    "G18" = "11"
    "D19" = "0.3319"
    "E19" = "-1.70%"
    "G19" = "11"
    "D20" = "7.051"
    "E20" = "-5.02%"
    "G20" = "11"
    "D21" = "0.1347"
    "E21" = "-3.02%"
    "G21" = "11"
    "D22" = "0.2762"
    "E22" = "8.44%"
    "G22" = "11"
    "D23" = "0.04135"
    "E23" = "-0.51%"
    "G23" = "11"
    "D24" = "0.001199"
    "E24" = "-3.61%"
    "G24" = "11"
    "D25" = "0.004070"
    "E25" = "-5.66%"
    "G25" = "11"
    "D26" = "0.0001302"
    "E26" = "8.37%"
    "G26" = "11"
    "G27" = "11"
    "G28" = "11"
    "G29" = "11"
    "G30" = "11"
    "G31" = "11"
    "G32" = "11"
    "G33" = "11"
    "G34" = "11"
    "G35" = "11"
    "G36" = "11"
    "G37" = "11"
    "D38" = "0.02342"
    "E38" = "-14.01%"
    "G38" = "11"
    "D39" = "0.05194"
    "E39" = "-6.53%"
    "G39" = "11"
    "D40" = "0.007936"
    "E40" = "-0.33%"
    "G40" = "11"
    "D41" = "0.1321"
    "E41" = "-5.44%"
    "G41" = "11"
    "D42" = "0.006458"
    "E42" = "-1.36%"
    "G42" = "11"
    "D43" = "0.001953"
    "E43" = "-6.47%"
    "G43" = "11"
    "D44" = "0.008737"
    "E44" = "5.44%"
    "G44" = "11"
    "D45" = "0.3329"
    "E45" = "-4.31%"
    "G45" = "11"
    "D46" = "0.00006520"
    "E46" = "-5.88%"
    "G46" = "11"
    "E47" = "0.12%"
    "G47" = "11"
    "E48" = "98.33%"
    "G48" = "11"
    "D49" = "0.003534"
    "E49" = "4.96%"
    "G49" = "11"
    "D50" = "0.00002104"
    "E50" = "0.12%"
    "G50" = "11"
    "D51" = "0.0002004"
    "E51" = "0.12%"
    "G51" = "11"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
}

Write-Host "Applied $($textUpdates.Count + $numericLookingUpdates.Count) cell updates"
